# Dispatch list: add a "UserID" column at the front of the sheet that
# carries a (possibly malformed) RedCap school id, so issues in that
# field can be checked while parsing RedCap data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remember where the existing mailto hyperlinks point, so we can
#     recreate them one column to the right after the insert (the engine
#     does not shift stored hyperlink refs when a column is inserted).
$oldLinks = @()
foreach ($h in $ws.Hyperlinks) {
    $oldLinks += [PSCustomObject]@{
        Addr   = $h.Range.Address()
        Target = $h.Address
    }
}

# Drop every existing hyperlink (one at a time from the live collection --
# iterating + deleting the same collection in one pass skips entries).
$linkCount = $ws.Hyperlinks.Count
for ($i = 0; $i -lt $linkCount; $i++) {
    foreach ($h in $ws.Hyperlinks) {
        $h.Delete()
        break
    }
}

# --- Insert a new column A; everything else shifts right by one.
$ws.Columns.Item(1).Insert()

# --- New header cell, bold like the other headers but its own style
#     (12pt bold instead of the 11pt bold used elsewhere).
$ws.Range("A1").Value = "UserID"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 12

# --- New per-row user/school identifiers.
$ws.Range("A2").Value = "AB9234"
$ws.Range("A3").Value = "AB100"
$ws.Range("A4").Value = "AB101"
$ws.Range("A5").Value = "another_school_180"

# --- Recreate the hyperlinks, shifted one column to the right of where
#     they used to be, and restore the cell's "Hyperlink" style
#     (Hyperlinks.Add always re-derives its own style xf, so we reset it
#     to the sheet's existing Hyperlink cell style explicitly).
function Shift-ColumnRight([string]$addr) {
    $addr = $addr -replace '\$', ''
    if ($addr -match '^([A-Z]+)([0-9]+)$') {
        $col = $matches[1]
        $row = $matches[2]
        $colNum = 0
        foreach ($ch in $col.ToCharArray()) {
            $colNum = $colNum * 26 + ([int][char]$ch - [int][char]'A' + 1)
        }
        $colNum = $colNum + 1
        $newCol = ""
        while ($colNum -gt 0) {
            $rem = ($colNum - 1) % 26
            $newCol = [string]([char]([int][char]'A' + $rem)) + $newCol
            $colNum = [int](($colNum - $rem) / 26)
        }
        return "$newCol$row"
    }
    return $addr
}

foreach ($l in $oldLinks) {
    $newAddr = Shift-ColumnRight($l.Addr)
    [void]$ws.Hyperlinks.Add($ws.Range($newAddr), $l.Target)
    $ws.Range($newAddr).Style = "Hyperlink"
}

# --- Leave the selection on A2, matching where the edit was made.
[void]$ws.Range("A2").Select()
